$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates (coin price / volume refresh) matching the commit diff.
$updates = [ordered]@{
    'D2' = '67.183.05'
    'E2' = '  -0.72%  '
    'D3' = '2.475.43'
    'E3' = '  -0.79%  '
    'E4' = '  -0.02%  '
    'D5' = '582.17'
    'E5' = '  -1.49%  '
    'D6' = '168.17'
    'E6' = '  -3.35%  '
    'E7' = '  +0.01%  '
    'D8' = '0.514'
    'E8' = '  -1.93%  '
    'D9' = '2.476.88'
    'E9' = '  -0.73%  '
    'D10' = '0.135'
    'E10' = '  -3.37%  '
    'E11' = '  -0.97%  '
    'E12' = '  -2.72%  '
    'D13' = '0.333'
    'E13' = '  -2.46%  '
    'D14' = '25.55'
    'E14' = '  -2.79%  '
    'D15' = '2.924.96'
    'E15' = '  -0.78%  '
    'D16' = '67.090.27'
    'E16' = '  -0.67%  '
    'D17' = '0.0000170'
    'E17' = '  -3.98%  '
    'D18' = '2.474.52'
    'E18' = '  -0.33%  '
    'D19' = '11.31'
    'E19' = '  -4.10%  '
    'D20' = '7.61'
    'E20' = '  -4.73%  '
    'D21' = '356.49'
    'E21' = '  -2.48%  '
    'D22' = '4.04'
    'E22' = '  -2.06%  '
    'E23' = '  +0.06%  '
    'D24' = '69.44'
    'E24' = '  -2.68%  '
    'D25' = '4.24'
    'E25' = '  -6.85%  '
    'D26' = '1.79'
    'E26' = '  -6.77%  '
    'D27' = '9.13'
    'E27' = '  -8.38%  '
    'D28' = '0.998'
    'E28' = '  -0.12%  '
    'E29' = '  -0.94%  '
    'D30' = '0.0₃0909'
    'E30' = '  -5.44%  '
    'D31' = '509.96'
    'E31' = '  -4.05%  '
    'D32' = '7.80'
    'E32' = '  -6.23%  '
    'E33' = '  -4.35%  '
    'E34' = '  -5.44%  '
    'E35' = '  +0.05%  '
    'E36' = '  -6.49%  '
    'D37' = '158.22'
    'E37' = '  +0.14%  '
    'B38' = 'EthereumClassic'
    'C38' = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
    'D38' = '18.50'
    'E38' = '  -1.07%  '
    'B39' = 'WhiteBITCoin'
    'C39' = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
    'D39' = '18.59'
    'E39' = '  -0.26%  '
    'D40' = '1.35'
    'E40' = '  -5.41%  '
    'E41' = '  +0.09%  '
    'D42' = '1.67'
    'E42' = '  -5.91%  '
    'B43' = 'RenderToken'
    'C43' = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
    'D43' = '4.81'
    'E43' = '  -5.80%  '
    'B44' = 'PolygonEcosystemToken'
    'C44' = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
    'D44' = '0.327'
    'E44' = '  -6.30%  '
    'D45' = '2.34'
    'E45' = '  -6.66%  '
    'D46' = '38.71'
    'E46' = '  -2.57%  '
    'D47' = '141.37'
    'E47' = '  -2.66%  '
    'D48' = '3.48'
    'E48' = '  -5.37%  '
    'D49' = '0.516'
    'E49' = '  -5.53%  '
    'D50' = '1.60'
    'E50' = '  -5.41%  '
    'D51' = '0.0₆0252'
    'E51' = '  -7.71%  '
}

# Every updated cell in this sheet is plain text (prices use "." as a thousands
# separator, not a decimal point, and percentages carry padding spaces). Any new
# value that merely *looks* numeric (e.g. "18.50", "0.0000170", "7.80") would be
# silently coerced by Excel into a Number cell, rounding/trimming it (18.50 -> 18.5,
# 0.0000170 -> 1.7E-05). Force Text format on those cells first so the exact
# literal string is preserved, matching the source XML (inline/shared string).
$textCells = @('D5', 'D6', 'D8', 'D10', 'D13', 'D14', 'D17', 'D19', 'D20', 'D21', 'D22', 'D24', 'D25', 'D26', 'D27', 'D28', 'D31', 'D32', 'D37', 'D38', 'D39', 'D40', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50')
foreach ($tc in $textCells) {
    $ws.Range($tc).NumberFormat = "@"
}

foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}

